$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column P width change (16th column): 8.71 -> 7.71 (matches column O's width)
$ws.Columns.Item(16).ColumnWidth = 6.83

# Updated precision estimation values for R11 edition 2.0 data
$ws.Range("P2").Value = 13.128
$ws.Range("P3").Value = 0.287
$ws.Range("U3").Value = 0.114
$ws.Range("P4").Value = 2.508
$ws.Range("U4").Value = 0.995
$ws.Range("P5").Value = 4.6
$ws.Range("U5").Value = 0.264
$ws.Range("P6").Value = 654.367
$ws.Range("U6").Value = 0.325
$ws.Range("P7").Value = 30.035
$ws.Range("U7").Value = 0.806
$ws.Range("P8").Value = 33.695
$ws.Range("U8").Value = 0.864
$ws.Range("P9").Value = 12.5
$ws.Range("U9").Value = 0.328
$ws.Range("P10").Value = 1958.3
$ws.Range("P11").Value = 1.871
$ws.Range("U11").Value = 0.993
$ws.Range("P12").Value = 4.966
$ws.Range("U12").Value = 0.933
$ws.Range("P13").Value = 1.939
$ws.Range("P14").Value = 2.037
$ws.Range("U14").Value = 0.983
$ws.Range("P15").Value = 7.464
$ws.Range("P16").Value = 2.144
$ws.Range("P17").Value = 5.271
$ws.Range("U17").Value = 0.972
$ws.Range("P18").Value = 5.335
$ws.Range("U18").Value = 0.975
$ws.Range("P19").Value = 4.979
$ws.Range("U19").Value = 0.972
$ws.Range("P20").Value = 4.001
$ws.Range("U20").Value = 0.895
$ws.Range("P21").Value = 4.411
$ws.Range("U21").Value = 0.877
$ws.Range("P22").Value = 2.622
$ws.Range("P23").Value = 5.679
$ws.Range("P24").Value = 5.095
$ws.Range("U24").Value = 0.997
$ws.Range("P25").Value = 5.149
$ws.Range("P26").Value = 4.765
$ws.Range("P27").Value = 5.054
$ws.Range("U27").Value = 0.969
$ws.Range("P28").Value = 4.505
$ws.Range("U28").Value = 0.984
$ws.Range("P29").Value = 5.122
$ws.Range("U29").Value = 0.958
$ws.Range("P30").Value = 4.058
$ws.Range("U30").Value = 0.97
$ws.Range("P31").Value = 5.106
$ws.Range("U31").Value = 0.992
$ws.Range("P32").Value = 7.167
$ws.Range("U32").Value = 0.993
$ws.Range("P33").Value = 4.323
$ws.Range("U33").Value = 0.94
$ws.Range("P34").Value = 5.463
$ws.Range("U34").Value = 0.982
$ws.Range("P35").Value = 6.505
$ws.Range("U35").Value = 0.992
$ws.Range("P36").Value = 3.671
$ws.Range("U36").Value = 0.985
$ws.Range("P37").Value = 4.524
$ws.Range("U37").Value = 0.981
$ws.Range("P38").Value = 3.631
$ws.Range("U38").Value = 0.982
$ws.Range("P39").Value = 4.777
$ws.Range("U39").Value = 0.933
$ws.Range("P40").Value = 5.669
$ws.Range("U40").Value = 0.898
$ws.Range("P41").Value = 0.268
$ws.Range("U41").Value = 0.996
$ws.Range("P42").Value = 0.899
$ws.Range("P43").Value = 0.347
$ws.Range("U43").Value = 0.675
$ws.Range("P44").Value = 0.44
$ws.Range("U44").Value = 0.977
$ws.Range("P45").Value = 0.109
$ws.Range("U45").Value = 0.998
$ws.Range("P46").Value = 0.024
$ws.Range("U46").Value = 0.456
$ws.Range("P47").Value = 0.945
$ws.Range("P50").Value = 0.006
$ws.Range("P51").Value = 0.01
$ws.Range("P52").Value = 0.015
$ws.Range("P53").Value = 0.077
$ws.Range("P56").Value = 0.922
$ws.Range("P58").Value = 0.016
$ws.Range("P63").Value = 0.155
$ws.Range("U63").Value = 0.995
$ws.Range("P64").Value = 0.771
$ws.Range("U64").Value = 0.902
$ws.Range("P65").Value = 0.864
$ws.Range("P66").Value = 0.467
$ws.Range("P67").Value = 0.313
$ws.Range("U67").Value = 0.883
$ws.Range("P68").Value = 0.869
$ws.Range("P69").Value = 0.341
$ws.Range("U69").Value = 0.432
$ws.Range("P70").Value = 0.63
$ws.Range("P71").Value = 0.12
$ws.Range("U71").Value = 0.365
$ws.Range("P73").Value = 0.107
$ws.Range("U73").Value = 0.256
$ws.Range("P74").Value = 0.258
$ws.Range("P75").Value = 0.097
$ws.Range("U75").Value = 0.257
$ws.Range("P76").Value = 0.736
$ws.Range("U76").Value = 0.991
$ws.Range("P77").Value = 0.062
$ws.Range("U77").Value = 0.903
